# Weekly update: insert a new "Ají" price-report row for Vega Central
# Mapocho de Santiago. This shifts the existing rows 348-365 down to
# 349-366 and fills the newly opened row 348 with this week's record.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 348, pushing rows 348:365 -> 349:366
$ws.Rows("348:348").Insert()

# Populate the new row 348 with the new weekly record
$ws.Range("A348").Value = 9
$ws.Range("B348").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C348").Value = "Metropolitana"
$ws.Range("D348").Value = 44931
$ws.Range("E348").Value = 13
$ws.Range("F348").Value = 100112021
$ws.Range("G348").Value = "Ají"
$ws.Range("H348").Value = "Inferno"
$ws.Range("I348").Value = "Primera"
$ws.Range("J348").Value = 90
$ws.Range("K348").Value = 12000
$ws.Range("L348").Value = 12000
$ws.Range("M348").Value = 12000
$ws.Range("N348").Value = "`$/caja 10 kilos"
$ws.Range("O348").Value = "Región de Arica y Parinacota"
$ws.Range("P348").Value = 1200
$ws.Range("Q348").Value = 10
$ws.Range("R348").Value = "Hortaliza"
